$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.430.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.647.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.87%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.644.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("E10").Value = "  +7.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.401"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.85"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000194"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +13.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.124.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.208.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.645.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "358.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.64%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.90%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000104"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("E29").Value = "  +1.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "524.20"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "162.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "165.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.87%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0607"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("E49").Value = "  +1.49%  "
$ws.Range("E50").Value = "  +3.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0978"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.08%  "
